# Repair for getting Properties from Object Repository
#
# TestSuite: mark every test case as Run Mode "Yes" / Result "PASSED"
# (previously most rows were "No" / blank), and move the active cell
# selection.
#
# TestCase2 .. TestCase9: fill the (previously blank) "Result" column G
# with "PASSED" for every data row. TestCase1 is left untouched.

$wb = $excel.ActiveWorkbook

# --- TestSuite sheet -------------------------------------------------
$tsSheet = $wb.Worksheets.Item("TestSuite")

$tsSheet.Range("D2").Value = "PASSED"
$tsSheet.Range("C3:C10").Value = "Yes"
$tsSheet.Range("D3:D10").Value = "PASSED"

$tsSheet.Range("B14").Select() | Out-Null

# --- TestCase2 .. TestCase9 sheets: fill Result column with PASSED ---
$resultRanges = @{
    "TestCase2" = "G2:G10"
    "TestCase3" = "G2:G10"
    "TestCase4" = "G2:G10"
    "TestCase5" = "G2:G9"
    "TestCase6" = "G2:G9"
    "TestCase7" = "G2:G7"
    "TestCase8" = "G2:G11"
    "TestCase9" = "G2:G8"
}

foreach ($sheetName in $resultRanges.Keys) {
    $sheet = $wb.Worksheets.Item($sheetName)
    $sheet.Range($resultRanges[$sheetName]).Value = "PASSED"
}
